$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 152 - shifts existing rows 152..192 down to 153..193
$ws.Rows.Item(152).Insert()

# Populate the newly inserted row 152 with the new record
$ws.Range('A152').Value = 4
$ws.Range('B152').Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range('C152').Value = 'Los Lagos'
$ws.Range('D152').Value = 44551
$ws.Range('E152').Value = 10
$ws.Range('F152').Value = 'Fruta'
$ws.Range('G152').Value = 100102
$ws.Range('H152').Value = 'Cítricos'
$ws.Range('I152').Value = 100102006
$ws.Range('J152').Value = 'Pomelo'
$ws.Range('K152').Value = 'Start Ruby'
$ws.Range('L152').Value = 'Primera'
$ws.Range('M152').Value = 200
$ws.Range('N152').Value = 10000
$ws.Range('O152').Value = 11000
$ws.Range('P152').Value = 10500
$ws.Range('Q152').Value = '$/caja 14 kilos empedrada'
$ws.Range('R152').Value = "Región de O'Higgins"
$ws.Range('S152').Value = 750
$ws.Range('T152').Value = 14
